$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all new cell values in the exact order they were authored
# so the shared-strings table is built up in the same sequence.
$ws.Range('B213').Value2 = 'BinaryOperationNode'
$ws.Range('A216').Value2 = 'operation'
$ws.Range('A217').Value2 = 'left'
$ws.Range('A218').Value2 = 'right'
$ws.Range('B223').Value2 = 'UnaryOperationNode'
$ws.Range('D216').Value2 = 'Операция в узле'
$ws.Range('D217').Value2 = 'Левое поддерево'
$ws.Range('D218').Value2 = 'Правое поддерево'
$ws.Range('D227').Value2 = 'Поддерево'
$ws.Range('B232').Value2 = 'AggregateOperationNode'
$ws.Range('C235').Value2 = 'AggregateOperation'
$ws.Range('D245').Value2 = 'Гистограмма в узле'
$ws.Range('B250').Value2 = 'SubhistogramNode'
$ws.Range('A253').Value2 = 'properties'
$ws.Range('A254').Value2 = 'originOpt'
$ws.Range('C254').Value2 = 'Option[Node[E]] '
$ws.Range('D253').Value2 = 'Подмножество'
$ws.Range('D254').Value2 = 'Поддерево узла'
$ws.Range('B259').Value2 = 'OperationInput'
$ws.Range('D262').Value2 = 'Вводимая операция'
$ws.Range('B264').Value2 = 'SubhistogramInput'
$ws.Range('D267').Value2 = 'Подмножество элементов'
$ws.Range('B269').Value2 = 'HistogramInput'
$ws.Range('D272').Value2 = 'Гистограмма'
$ws.Range('B274').Value2 = 'Parser'
$ws.Range('A277').Value2 = 'parse'
$ws.Range('C277').Value2 = 'Option[Stack[Input[E]]]'
$ws.Range('D277').Value2 = 'query: String, implicit aliasToInput: Map[String, Input[E]]'
$ws.Range('A278').Value2 = 'getLexems'
$ws.Range('D278').Value2 = 'query: String, acc: Stack[Input[E]], aliasToInput: Map[String, Input[E]]'
$ws.Range('A279').Value2 = 'toPolishNotation'
$ws.Range('C279').Value2 = 'Stack[Input[E]] '
$ws.Range('D279').Value2 = 'query: Stack[Input[E]], resultAcc: Stack[Input[E]], operandsAcc: Stack[Input[E]]'
$ws.Range('E279').Value2 = 'Преобразует последовательность входов в польскую инверсную последовательность'
$ws.Range('E278').Value2 = 'Преобразует строку в последовательность входных конструкций'
$ws.Range('E277').Value2 = 'Преобразует строку в последовательность входных лексем в польской нотации'
$ws.Range('B281').Value2 = 'TreeExecutor'
$ws.Range('A284').Value2 = 'execute'
$ws.Range('C284').Value2 = 'Either[Histogram[E], Double]'
$ws.Range('D284').Value2 = 'tree: Node[E]'
$ws.Range('E284').Value2 = 'Исполняет вычисления на АСТ'
$ws.Range('B286').Value2 = 'Query'
$ws.Range('A289').Value2 = 'root'
$ws.Range('D289').Value2 = 'АСТ'
$ws.Range('A294').Value2 = 'parseStack'
$ws.Range('A296').Value2 = 'fromString'
$ws.Range('A290').Value2 = 'standardAliases'
$ws.Range('C290').Value2 = 'Map[String, Input[E]]'
$ws.Range('D290').Value2 = 'Стандартные алиасы входных конструкций'
$ws.Range('E293').Value2 = 'Исполняет запрос'
$ws.Range('E294').Value2 = 'Преобразует польскую запись в АСТ входных конструкций'
$ws.Range('D294').Value2 = 'operationsStack: Stack[Input[E]]'
$ws.Range('C295').Value2 = 'Query[E]'
$ws.Range('D295').Value2 = 'Stack[Input[E]]'
$ws.Range('E295').Value2 = 'Формирует запрос по последовательности входных конструкций'
$ws.Range('D296').Value2 = 'query: String, aliasToInput: Map[String, Input[E]]'
$ws.Range('E296').Value2 = 'Формирует запрос по строковому представлению'
$ws.Range('A213').Value2 = 'Класс'
$ws.Range('A214').Value2 = 'Поля'
$ws.Range('A215').Value2 = 'Имя'
$ws.Range('B215').Value2 = 'Модификатор доступа'
$ws.Range('C215').Value2 = 'Тип'
$ws.Range('D215').Value2 = 'Назначение'
$ws.Range('B216').Value2 = '-'
$ws.Range('C216').Value2 = 'HistogramBinaryOperation'
$ws.Range('B217').Value2 = '-'
$ws.Range('C217').Value2 = 'Node[E]'
$ws.Range('B218').Value2 = '-'
$ws.Range('C218').Value2 = 'Node[E]'
$ws.Range('A219').Value2 = 'Методы'
$ws.Range('A220').Value2 = 'Имя'
$ws.Range('B220').Value2 = 'Модификатор доступа'
$ws.Range('C220').Value2 = 'Тип'
$ws.Range('D220').Value2 = 'Аргументы'
$ws.Range('E220').Value2 = 'Назначение'
$ws.Range('A221').Value2 = 'map'
$ws.Range('B221').Value2 = '-'
$ws.Range('C221').Value2 = 'Node[E]'
$ws.Range('D221').Value2 = 'f: Node[E] => Node[E]'
$ws.Range('E221').Value2 = 'Преобразует узел дерева'
$ws.Range('A223').Value2 = 'Класс'
$ws.Range('A224').Value2 = 'Поля'
$ws.Range('A225').Value2 = 'Имя'
$ws.Range('B225').Value2 = 'Модификатор доступа'
$ws.Range('C225').Value2 = 'Тип'
$ws.Range('D225').Value2 = 'Назначение'
$ws.Range('A226').Value2 = 'operation'
$ws.Range('B226').Value2 = '-'
$ws.Range('C226').Value2 = 'HistogramUnaryOperation'
$ws.Range('D226').Value2 = 'Операция в узле'
$ws.Range('A227').Value2 = 'histogram'
$ws.Range('B227').Value2 = '-'
$ws.Range('C227').Value2 = 'Node[E]'
$ws.Range('A228').Value2 = 'Методы'
$ws.Range('A229').Value2 = 'Имя'
$ws.Range('B229').Value2 = 'Модификатор доступа'
$ws.Range('C229').Value2 = 'Тип'
$ws.Range('D229').Value2 = 'Аргументы'
$ws.Range('E229').Value2 = 'Назначение'
$ws.Range('A230').Value2 = 'map'
$ws.Range('B230').Value2 = '-'
$ws.Range('C230').Value2 = 'Node[E]'
$ws.Range('D230').Value2 = 'f: Node[E] => Node[E]'
$ws.Range('E230').Value2 = 'Преобразует узел дерева'
$ws.Range('A232').Value2 = 'Класс'
$ws.Range('A233').Value2 = 'Поля'
$ws.Range('A234').Value2 = 'Имя'
$ws.Range('B234').Value2 = 'Модификатор доступа'
$ws.Range('C234').Value2 = 'Тип'
$ws.Range('D234').Value2 = 'Назначение'
$ws.Range('A235').Value2 = 'operation'
$ws.Range('B235').Value2 = '-'
$ws.Range('D235').Value2 = 'Операция в узле'
$ws.Range('A236').Value2 = 'left'
$ws.Range('B236').Value2 = '-'
$ws.Range('C236').Value2 = 'Node[E]'
$ws.Range('D236').Value2 = 'Левое поддерево'
$ws.Range('A237').Value2 = 'right'
$ws.Range('B237').Value2 = '-'
$ws.Range('C237').Value2 = 'Node[E]'
$ws.Range('D237').Value2 = 'Правое поддерево'
$ws.Range('A238').Value2 = 'Методы'
$ws.Range('A239').Value2 = 'Имя'
$ws.Range('B239').Value2 = 'Модификатор доступа'
$ws.Range('C239').Value2 = 'Тип'
$ws.Range('D239').Value2 = 'Аргументы'
$ws.Range('E239').Value2 = 'Назначение'
$ws.Range('A240').Value2 = 'map'
$ws.Range('B240').Value2 = '-'
$ws.Range('C240').Value2 = 'Node[E]'
$ws.Range('D240').Value2 = 'f: Node[E] => Node[E]'
$ws.Range('E240').Value2 = 'Преобразует узел дерева'
$ws.Range('A242').Value2 = 'Класс'
$ws.Range('B242').Value2 = 'AggregateOperationNode'
$ws.Range('A243').Value2 = 'Поля'
$ws.Range('A244').Value2 = 'Имя'
$ws.Range('B244').Value2 = 'Модификатор доступа'
$ws.Range('C244').Value2 = 'Тип'
$ws.Range('D244').Value2 = 'Назначение'
$ws.Range('A245').Value2 = 'histogram'
$ws.Range('B245').Value2 = '-'
$ws.Range('C245').Value2 = 'AggregateOperation'
$ws.Range('A246').Value2 = 'Методы'
$ws.Range('A247').Value2 = 'Имя'
$ws.Range('B247').Value2 = 'Модификатор доступа'
$ws.Range('C247').Value2 = 'Тип'
$ws.Range('D247').Value2 = 'Аргументы'
$ws.Range('E247').Value2 = 'Назначение'
$ws.Range('A248').Value2 = 'map'
$ws.Range('B248').Value2 = '-'
$ws.Range('C248').Value2 = 'Node[E]'
$ws.Range('D248').Value2 = 'f: Node[E] => Node[E]'
$ws.Range('E248').Value2 = 'Преобразует узел дерева'
$ws.Range('A250').Value2 = 'Класс'
$ws.Range('A251').Value2 = 'Поля'
$ws.Range('A252').Value2 = 'Имя'
$ws.Range('B252').Value2 = 'Модификатор доступа'
$ws.Range('C252').Value2 = 'Тип'
$ws.Range('D252').Value2 = 'Назначение'
$ws.Range('B253').Value2 = '-'
$ws.Range('C253').Value2 = 'ElementsUniverse[E]'
$ws.Range('B254').Value2 = '-'
$ws.Range('A255').Value2 = 'Методы'
$ws.Range('A256').Value2 = 'Имя'
$ws.Range('B256').Value2 = 'Модификатор доступа'
$ws.Range('C256').Value2 = 'Тип'
$ws.Range('D256').Value2 = 'Аргументы'
$ws.Range('E256').Value2 = 'Назначение'
$ws.Range('A257').Value2 = 'map'
$ws.Range('B257').Value2 = '-'
$ws.Range('C257').Value2 = 'Node[E]'
$ws.Range('D257').Value2 = 'f: Node[E] => Node[E]'
$ws.Range('E257').Value2 = 'Преобразует узел дерева'
$ws.Range('A259').Value2 = 'Класс'
$ws.Range('A260').Value2 = 'Поля'
$ws.Range('A261').Value2 = 'Имя'
$ws.Range('B261').Value2 = 'Модификатор доступа'
$ws.Range('C261').Value2 = 'Тип'
$ws.Range('D261').Value2 = 'Назначение'
$ws.Range('A262').Value2 = 'operation'
$ws.Range('B262').Value2 = '-'
$ws.Range('C262').Value2 = 'Operation'
$ws.Range('A264').Value2 = 'Класс'
$ws.Range('A265').Value2 = 'Поля'
$ws.Range('A266').Value2 = 'Имя'
$ws.Range('B266').Value2 = 'Модификатор доступа'
$ws.Range('C266').Value2 = 'Тип'
$ws.Range('D266').Value2 = 'Назначение'
$ws.Range('A267').Value2 = 'properties'
$ws.Range('B267').Value2 = '-'
$ws.Range('C267').Value2 = 'ElementsUniverse[E]'
$ws.Range('A269').Value2 = 'Класс'
$ws.Range('A270').Value2 = 'Поля'
$ws.Range('A271').Value2 = 'Имя'
$ws.Range('B271').Value2 = 'Модификатор доступа'
$ws.Range('C271').Value2 = 'Тип'
$ws.Range('D271').Value2 = 'Назначение'
$ws.Range('A272').Value2 = 'histogram'
$ws.Range('B272').Value2 = '-'
$ws.Range('C272').Value2 = 'Histogram[E]'
$ws.Range('A274').Value2 = 'Класс'
$ws.Range('A275').Value2 = 'Методы'
$ws.Range('A276').Value2 = 'Имя'
$ws.Range('B276').Value2 = 'Модификатор доступа'
$ws.Range('C276').Value2 = 'Тип'
$ws.Range('D276').Value2 = 'Аргументы'
$ws.Range('E276').Value2 = 'Назначение'
$ws.Range('B277').Value2 = '-'
$ws.Range('B278').Value2 = '-'
$ws.Range('C278').Value2 = 'Option[Stack[Input[E]]]'
$ws.Range('B279').Value2 = '-'
$ws.Range('A281').Value2 = 'Класс'
$ws.Range('A282').Value2 = 'Методы'
$ws.Range('A283').Value2 = 'Имя'
$ws.Range('B283').Value2 = 'Модификатор доступа'
$ws.Range('C283').Value2 = 'Тип'
$ws.Range('D283').Value2 = 'Аргументы'
$ws.Range('E283').Value2 = 'Назначение'
$ws.Range('B284').Value2 = '-'
$ws.Range('A286').Value2 = 'Класс'
$ws.Range('A287').Value2 = 'Поля'
$ws.Range('A288').Value2 = 'Имя'
$ws.Range('B288').Value2 = 'Модификатор доступа'
$ws.Range('C288').Value2 = 'Тип'
$ws.Range('D288').Value2 = 'Назначение'
$ws.Range('B289').Value2 = '-'
$ws.Range('C289').Value2 = 'Node[E]'
$ws.Range('B290').Value2 = '-'
$ws.Range('A291').Value2 = 'Методы'
$ws.Range('A292').Value2 = 'Имя'
$ws.Range('B292').Value2 = 'Модификатор доступа'
$ws.Range('C292').Value2 = 'Тип'
$ws.Range('D292').Value2 = 'Аргументы'
$ws.Range('E292').Value2 = 'Назначение'
$ws.Range('A293').Value2 = 'execute'
$ws.Range('B293').Value2 = '-'
$ws.Range('C293').Value2 = 'Either[Histogram[E], Double]'
$ws.Range('D293').Value2 = 'histogram: Histogram[E]'
$ws.Range('B294').Value2 = '-'
$ws.Range('C294').Value2 = 'Node[E]'
$ws.Range('A295').Value2 = 'apply'
$ws.Range('B295').Value2 = '-'
$ws.Range('B296').Value2 = '-'
$ws.Range('C296').Value2 = 'Query[E]'

# Apply wrap-text formatting to the cells that need it
$ws.Range('D278').WrapText = $true
$ws.Range('D285').WrapText = $true
$ws.Range('D294').WrapText = $true

# Set explicit row heights for the wrapped rows
$ws.Rows.Item(278).RowHeight = 48
$ws.Rows.Item(294).RowHeight = 32

# Update the sheet selection to match the final cursor position
$ws.Range('E297').Select()
